$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D shifts to F, etc.)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formatting/styles from the (now-shifted) old column D - now column F -
# into the two newly inserted columns D:E, restricted to the contiguous row blocks
# that actually carry data (this avoids manufacturing stray blank cells on label-only
# rows such as 5, 6, 37 and 79).
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)

$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)

$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Populate the two new columns with the new quarter figures.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 1013200
$ws.Range("E8").Value = 888400
$ws.Range("D9").Value = 889700
$ws.Range("E9").Value = 788000
$ws.Range("D10").Value = 123500
$ws.Range("E10").Value = 100400
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 921400
$ws.Range("E17").Value = 817000
$ws.Range("D18").Value = 91800
$ws.Range("E18").Value = 71400
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 99300
$ws.Range("E21").Value = 78300
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 91800
$ws.Range("E23").Value = 71400
$ws.Range("D24").Value = 19000
$ws.Range("E24").Value = 17300
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 72700
$ws.Range("E26").Value = 54100
$ws.Range("D27").Value = 72700
$ws.Range("E27").Value = 54100
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 2700
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 75500
$ws.Range("E33").Value = 54100
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 75500
$ws.Range("E35").Value = 54100
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 311500
$ws.Range("E41").Value = 205800
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 77300
$ws.Range("E43").Value = 79600
$ws.Range("D44").Value = 2742600
$ws.Range("E44").Value = 2887300
$ws.Range("D45").Value = 84200
$ws.Range("E45").Value = 82800
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 17500
$ws.Range("E47").Value = 16300
$ws.Range("D48").Value = 54600
$ws.Range("E48").Value = 53400
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 26500
$ws.Range("E52").Value = 36700
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 3365500
$ws.Range("E54").Value = 3448300
$ws.Range("D57").Value = 128200
$ws.Range("E57").Value = 156800
$ws.Range("D58").Value = 14800
$ws.Range("E58").Value = 16700
$ws.Range("D59").Value = 206500
$ws.Range("E59").Value = 234600
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 1295300
$ws.Range("E61").Value = 1295100
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1644700
$ws.Range("E66").Value = 1735800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 1218600
$ws.Range("E72").Value = 1143100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1720800
$ws.Range("E76").Value = 1712500
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 75500
$ws.Range("E81").Value = 54100
$ws.Range("D83").Value = 7500
$ws.Range("E83").Value = 6900
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 188500
$ws.Range("E89").Value = 84100
$ws.Range("D91").Value = -9700
$ws.Range("E91").Value = -8000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -9900
$ws.Range("E94").Value = -7600
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -72900
$ws.Range("E100").Value = -40200
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 105700
$ws.Range("E102").Value = 36300
